$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 ---
$ws.Range("B2").Value = "lacy"
$ws.Range("C2").Value = 45546
$ws.Range("E2").Value = "Macho"
$ws.Range("G2").Value = 1

# --- Add row 3 ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "lacy"
$ws.Range("C3").Value = 45546
$ws.Range("C3").NumberFormat = "YYYY-MM-DD"
$ws.Range("D3").Value = "dorper"
$ws.Range("E3").Value = "Macho"

# --- Add row 4 ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "lacy"
$ws.Range("C4").Value = 45546
$ws.Range("C4").NumberFormat = "YYYY-MM-DD"
$ws.Range("D4").Value = "dorper"
$ws.Range("E4").Value = "Macho"
